$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stage a copy of every distinct fill/font style that already exists in
#    rows 6-14 into an out-of-the-way scratch column (Z), so it can be
#    reapplied after the rows are reshuffled without clobbering a style
#    before it has been captured.
#      Z1 <- style 1  (red fill)                 source B6
#      Z2 <- style 2  (yellow fill)               source C11
#      Z3 <- style 3  (blue fill)                 source B9
#      Z4 <- style 4  (green fill)                source C6
#      Z5 <- style 5  (orange fill)                source B7
#      Z6 <- style 7  (yellow fill + yellow font) source C9
#      Z7 <- style 8  (red fill + plain font)     source B10
#      Z8 <- style 9  (yellow fill + valign)      source C10
#      Z9 <- style 10 (green fill + yellow font)  source C14
# ---------------------------------------------------------------------------
$ws.Range("B6").Copy($ws.Range("Z1"))
$ws.Range("C11").Copy($ws.Range("Z2"))
$ws.Range("B9").Copy($ws.Range("Z3"))
$ws.Range("C6").Copy($ws.Range("Z4"))
$ws.Range("B7").Copy($ws.Range("Z5"))
$ws.Range("C9").Copy($ws.Range("Z6"))
$ws.Range("B10").Copy($ws.Range("Z7"))
$ws.Range("C10").Copy($ws.Range("Z8"))
$ws.Range("C14").Copy($ws.Range("Z9"))

# ---------------------------------------------------------------------------
# 2) Add the three new tasks (new shared strings).
# ---------------------------------------------------------------------------
$ws.Range("A20000").Value = "Мультиплеер"
$ws.Range("A20001").Value = "Кооператив"
$ws.Range("A20002").Value = "УДАЛЕНИП КАСТОМНЫХ СПРАЙТОВ"

# ---------------------------------------------------------------------------
# 3) Rebuild rows 6-17 in the new (color-sorted) order. Clear each row
#    first so stale cells (e.g. an old C column) don't linger, then set
#    the task name and reapply the correct fill style from the scratch
#    column staged above.
# ---------------------------------------------------------------------------

# Row 6: ИИ
$ws.Range("A6:D6").Clear()
$ws.Range("A6").Value = "ИИ"
$ws.Range("Z7").Copy($ws.Range("B6"))
$ws.Range("Z8").Copy($ws.Range("C6"))

# Row 7: Маски ударов
$ws.Range("A7:D7").Clear()
$ws.Range("A7").Value = "Маски ударов"
$ws.Range("Z1").Copy($ws.Range("B7"))
$ws.Range("Z4").Copy($ws.Range("C7"))

# Row 8: Хитбоксы и кикбоксы
$ws.Range("A8:D8").Clear()
$ws.Range("A8").Value = "Хитбоксы и кикбоксы"
$ws.Range("Z1").Copy($ws.Range("B8"))
$ws.Range("Z9").Copy($ws.Range("C8"))

# Row 9: Больше оружия (unchanged position/content, rewritten for consistency)
$ws.Range("A9:D9").Clear()
$ws.Range("A9").Value = "Больше оружия"
$ws.Range("Z3").Copy($ws.Range("B9"))
$ws.Range("Z6").Copy($ws.Range("C9"))

# Row 10: Сундуки
$ws.Range("A10:D10").Clear()
$ws.Range("A10").Value = "Сундуки"
$ws.Range("Z3").Copy($ws.Range("B10"))
$ws.Range("Z4").Copy($ws.Range("C10"))

# Row 11: Покупки
$ws.Range("A11:D11").Clear()
$ws.Range("A11").Value = "Покупки"
$ws.Range("Z3").Copy($ws.Range("B11"))
$ws.Range("Z4").Copy($ws.Range("C11"))

# Row 12: Мультиплеер (new task, no C column)
$ws.Range("A12:D12").Clear()
$ws.Range("A12").Value = "Мультиплеер"
$ws.Range("Z3").Copy($ws.Range("B12"))

# Row 13: Кооператив (new task, no C column)
$ws.Range("A13:D13").Clear()
$ws.Range("A13").Value = "Кооператив"
$ws.Range("Z3").Copy($ws.Range("B13"))

# Row 14: Генерация уровней
$ws.Range("A14:D14").Clear()
$ws.Range("A14").Value = "Генерация уровней"
$ws.Range("Z5").Copy($ws.Range("B14"))
$ws.Range("Z2").Copy($ws.Range("C14"))

# Row 15: Двойной удар афины
$ws.Range("A15:D15").Clear()
$ws.Range("A15").Value = "Двойной удар афины"
$ws.Range("Z5").Copy($ws.Range("B15"))
$ws.Range("Z4").Copy($ws.Range("C15"))

# Row 16: 3 удар крит у топора
$ws.Range("A16:D16").Clear()
$ws.Range("A16").Value = "3 удар крит у топора"
$ws.Range("Z5").Copy($ws.Range("B16"))
$ws.Range("Z4").Copy($ws.Range("C16"))

# Row 17: УДАЛЕНИП КАСТОМНЫХ СПРАЙТОВ (new task, A column only)
$ws.Range("A17:D17").Clear()
$ws.Range("A17").Value = "УДАЛЕНИП КАСТОМНЫХ СПРАЙТОВ"

# ---------------------------------------------------------------------------
# 4) Drop the scratch column and the temporary helper cells.
# ---------------------------------------------------------------------------
$ws.Range("Z1:Z9").Clear()
$ws.Range("A20000:A20002").Clear()

# ---------------------------------------------------------------------------
# 5) Recreate the two differential (highlight) formats used by the color
#    sort, then clear the conditional-formatting rule itself so only the
#    dxf definitions remain (matches the committed workbook, which has
#    dxfs but no live conditionalFormatting entries).
# ---------------------------------------------------------------------------
$fcYellow = $ws.Range("B6:B16").FormatConditions.Add(2, 5, "1")
$fcYellow.Interior.Color = 65535
$fcRed = $ws.Range("C6:C16").FormatConditions.Add(2, 5, "1")
$fcRed.Interior.Color = 255
$ws.Range("B6:B16").FormatConditions.Delete()
$ws.Range("C6:C16").FormatConditions.Delete()

# ---------------------------------------------------------------------------
# 6) Record the sort state for the table (sorted by cell color).
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B6:B16"), 3, 1, $null, 0)
$ws.Sort.SortFields.Add($ws.Range("C6:C16"), 3, 1, $null, 0)
$ws.Sort.SetRange($ws.Range("A6:C16"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 7) Update the view: drop the old scroll position and select D11.
# ---------------------------------------------------------------------------
$ws.Range("D11").Select()
